$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(4101, 4544, 4544, 4544, 4689, 4689, 4932, 4932, 4932, 4932, 4932, 4932, 4932, 4932)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
